# Inserts a new daily price record at row 140 of the single data sheet,
# shifting every existing row from 140..274 down by one (to 141..275),
# matching the "weekly -> extra daily observation" update described in the
# commit message ("Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 140..274 down to 141..275 and free up row 140 for the new record.
$ws.Rows(140).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A140").Value = 8
$ws.Range("B140").Value = "Terminal La Palmera de La Serena"
$ws.Range("C140").Value = "Coquimbo"
$ws.Range("D140").Value = 44587
$ws.Range("E140").Value = 4
$ws.Range("F140").Value = 100114013
$ws.Range("G140").Value = "Zanahoria"
$ws.Range("H140").Value = "Sin especificar"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 800
$ws.Range("K140").Value = 5500
$ws.Range("L140").Value = 6000
$ws.Range("M140").Value = 5750
$ws.Range("N140").Value = "$/saco 20 kilos"
$ws.Range("O140").Value = "Provincia del Elquí"
$ws.Range("P140").Value = 288
$ws.Range("Q140").Value = 20
$ws.Range("R140").Value = "Hortaliza"
